$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Payouts")

# Remove the "Working Name" (G) and "Real Name" (H) columns entirely.
# This shifts "Notes" (old column I) left into column G.
$ws.Range("G1:H1").EntireColumn.Delete()
